# Automatische test-sync: 2025-06-30 20:01:50
# Adds a new log row (#11) to the "Logs" sheet and updates the
# "Dashboard" summary count for the "Retour / Terugbetaling" category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append new row 11 with the test mail data -----------------------------
$ws.Range("A11").Value = "Mijn product is beschadigd geleverd."
$ws.Range("B11").Value = "mailmind.test@zohomail.eu"
$ws.Range("C11").Value = "Testmail #11: Mijn product is beschadigd geleverd."
$ws.Range("D11").Value = "Retour / Terugbetaling"
$ws.Range("E11").Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om u beter van dienst te kunnen zijn, hebben wij wat meer informatie nodig over de beschadiging. Kunt u alstublieft een foto van het beschadigde product meesturen? Dit helpt ons om het probleem beter te begrijpen en een passende oplossing voor u te vinden.`nMet vriendelijke groet,`n[E-mailassistent] van [Bedrijfsnaam]"
$ws.Range("F11").Value = "2025-06-30 20:01:14"
$ws.Range("G11").Value = "Ja"
$ws.Range("H11").Value = "Nee"
$ws.Range("I11").Value = "Ja"
$ws.Range("J11").Value = "Nee"

# Writing the multi-line answer into E11 makes Excel apply a custom row
# height; restore the standard (auto-fit) height so the row matches the
# other, non-customized rows.
$ws.Rows.Item(11).AutoFit()

# --- Extend the conditional formatting ranges to include the new row -------
$ws.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D11"))
$ws.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G11"))
$ws.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H11"))
$ws.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I11"))
$ws.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J11"))

# --- Update the Dashboard summary count for "Retour / Terugbetaling" -------
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B2").Value = 4
